# Weekly update: insert the newest week's two "Coliflor" price observations
# (Primera / Segunda) at the top of the data block, pushing the existing
# history down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 499:500 - existing rows 499:592 shift to 501:594,
# inheriting formatting (e.g. the date style on column D) from the row above,
# exactly like an interactive Excel "Insert Copied Cells" / row insert.
$ws.Rows("499:500").Insert()

# New row 499: Primera
$ws.Range("A499").Value = 7
$ws.Range("B499").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C499").Value = "Ñuble"
$ws.Range("D499").Value = 45211
$ws.Range("E499").Value = 16
$ws.Range("F499").Value = 100112008
$ws.Range("G499").Value = "Coliflor"
$ws.Range("H499").Value = "Sin especificar"
$ws.Range("I499").Value = "Primera"
$ws.Range("J499").Value = 500
$ws.Range("K499").Value = 1200
$ws.Range("L499").Value = 1200
$ws.Range("M499").Value = 1200
$ws.Range("N499").Value = "`$/unidad"
$ws.Range("O499").Value = "Región del Maule"
$ws.Range("P499").Value = 1200
$ws.Range("Q499").Value = 1
$ws.Range("R499").Value = "Hortaliza"

# New row 500: Segunda
$ws.Range("A500").Value = 7
$ws.Range("B500").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C500").Value = "Ñuble"
$ws.Range("D500").Value = 45211
$ws.Range("E500").Value = 16
$ws.Range("F500").Value = 100112008
$ws.Range("G500").Value = "Coliflor"
$ws.Range("H500").Value = "Sin especificar"
$ws.Range("I500").Value = "Segunda"
$ws.Range("J500").Value = 400
$ws.Range("K500").Value = 1000
$ws.Range("L500").Value = 1000
$ws.Range("M500").Value = 1000
$ws.Range("N500").Value = "`$/unidad"
$ws.Range("O500").Value = "Región del Maule"
$ws.Range("P500").Value = 1000
$ws.Range("Q500").Value = 1
$ws.Range("R500").Value = "Hortaliza"
